# Insert a new weekly price record for "Albahaca" right after the existing
# row 94 (i.e. at row 95), pushing the former rows 95-123 down to 96-124.
# The new row carries the same attributes as the record that used to sit
# at row 95, except for an updated date (Fecha) and volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 95..123 down to 96..124, leaving a blank row 95 behind.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new record.
$ws.Cells.Item(95, 1).Value  = 8
$ws.Cells.Item(95, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(95, 3).Value  = "Coquimbo"
$ws.Cells.Item(95, 4).Value  = 44839
$ws.Cells.Item(95, 5).Value  = 4
$ws.Cells.Item(95, 6).Value  = 100112052
$ws.Cells.Item(95, 7).Value  = "Albahaca"
$ws.Cells.Item(95, 8).Value  = "Sin especificar"
$ws.Cells.Item(95, 9).Value  = "Primera"
$ws.Cells.Item(95, 10).Value = 1400
$ws.Cells.Item(95, 11).Value = 4000
$ws.Cells.Item(95, 12).Value = 4500
$ws.Cells.Item(95, 13).Value = 4250
$ws.Cells.Item(95, 14).Value = "$/paquete"
$ws.Cells.Item(95, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(95, 16).Value = 4250
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = "Hortaliza"
